$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.667.18"
$ws.Range("E2").Value = "'  +2.24%  "
$ws.Range("D3").Value = "'3.103.68"
$ws.Range("E3").Value = "'  +0.77%  "
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("D5").Value = "'527.47"
$ws.Range("E5").Value = "'  +2.36%  "
$ws.Range("D6").Value = "'143.85"
$ws.Range("E6").Value = "'  +1.84%  "
$ws.Range("E7").Value = "'  -0.06%  "
$ws.Range("E8").Value = "'  +1.67%  "
$ws.Range("D9").Value = "'7.36"
$ws.Range("E9").Value = "'  +1.32%  "
$ws.Range("E10").Value = "'  +0.80%  "
$ws.Range("D11").Value = "'0.383"
$ws.Range("E11").Value = "'  +2.82%  "
$ws.Range("D12").Value = "'3.638.59"
$ws.Range("E12").Value = "'  +0.77%  "
$ws.Range("E13").Value = "'  +1.05%  "
$ws.Range("D14").Value = "'26.92"
$ws.Range("E14").Value = "'  +5.46%  "
$ws.Range("E15").Value = "'  +2.31%  "
$ws.Range("D16").Value = "'58.659.57"
$ws.Range("E16").Value = "'  +2.03%  "
$ws.Range("D17").Value = "'3.101.54"
$ws.Range("E17").Value = "'  +0.79%  "
$ws.Range("E18").Value = "'  +0.66%  "
$ws.Range("D19").Value = "'12.96"
$ws.Range("E19").Value = "'  -1.14%  "
$ws.Range("D21").Value = "'342.22"
$ws.Range("E21").Value = "'  +2.41%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "'  -0.28%  "
$ws.Range("D23").Value = "'0.506"
$ws.Range("E23").Value = "'  +1.15%  "
$ws.Range("D24").Value = "'65.99"
$ws.Range("E24").Value = "'  +0.01%  "
$ws.Range("D25").Value = "'0.170"
$ws.Range("E25").Value = "'  +0.47%  "
$ws.Range("E26").Value = "'  -0.04%  "
$ws.Range("D27").Value = "'0.0₃0919"
$ws.Range("E27").Value = "'  +1.31%  "
$ws.Range("E28").Value = "'  +4.66%  "
$ws.Range("E29").Value = "'  +2.42%  "
$ws.Range("E30").Value = "'  +3.09%  "
$ws.Range("E31").Value = "'  +3.59%  "
$ws.Range("D32").Value = "'20.99"
$ws.Range("E32").Value = "'  +0.75%  "
$ws.Range("D33").Value = "'153.94"
$ws.Range("E33").Value = "'  -0.35%  "
$ws.Range("D34").Value = "'4.66"
$ws.Range("E34").Value = "'  +2.89%  "
$ws.Range("D35").Value = "'6.08"
$ws.Range("E35").Value = "'  +2.88%  "
$ws.Range("D36").Value = "'27.26"
$ws.Range("E36").Value = "'  -2.36%  "
$ws.Range("E37").Value = "'  +3.86%  "
$ws.Range("E38").Value = "'  +0.50%  "
$ws.Range("D39").Value = "'3.145.61"
$ws.Range("E39").Value = "'  +0.77%  "
$ws.Range("E40").Value = "'  +1.53%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.51"
$ws.Range("E41").Value = "'  +9.42%  "
$ws.Range("D42").Value = "'36.88"
$ws.Range("E42").Value = "'  +0.26%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.677"
$ws.Range("E43").Value = "'  +1.06%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "'  -0.03%  "
$ws.Range("D45").Value = "'2.297.43"
$ws.Range("E45").Value = "'  -0.10%  "
$ws.Range("D46").Value = "'0.0258"
$ws.Range("E46").Value = "'  +1.70%  "
$ws.Range("D47").Value = "'21.05"
$ws.Range("E47").Value = "'  +4.86%  "
$ws.Range("D48").Value = "'0.974"
$ws.Range("E48").Value = "'  +3.55%  "
$ws.Range("E49").Value = "'  +2.11%  "
$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").Value = "'270.56"
$ws.Range("E50").Value = "'  +6.91%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").Value = "'0.753"
$ws.Range("E51").Value = "'  +9.59%  "
